# Trade #51 closed at 2026-02-17 15:42:43 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet updates
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.54   # Current Capital
$summary.Range("B4").Value = 0.54      # Total P&L $
$summary.Range("B5").Value = 0.21      # Total P&L %
$summary.Range("B6").Value = 51        # Total Trades
$summary.Range("B7").Value = 15        # Winning Trades
$summary.Range("B9").Value = 29.41     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet updates (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.54
$status.Range("D4").Value = 51
$status.Range("E4").Value = 0.54
$status.Range("F4").Value = 0.54
$status.Range("G4").Value = 29.41

# ---------------------------------------------------------------------
# Helper to append the new trade (#51) row to a trades sheet
# ---------------------------------------------------------------------
function Add-TradeRow($sheet) {
    $row = 52

    $sheet.Cells.Item($row, 1).Value = 51

    # Date/Time columns hold plain text like "2026-02-17" / "15:42:37" in
    # this workbook (not real Excel date/time values). Force the cell to
    # Text format first so Excel's auto-detection doesn't turn the date
    # string into a date serial number.
    $sheet.Cells.Item($row, 2).NumberFormat = "@"
    $sheet.Cells.Item($row, 2).Value = "2026-02-17"

    $sheet.Cells.Item($row, 3).Value = "15:42:37"
    $sheet.Cells.Item($row, 4).Value = "MarketMaking"
    $sheet.Cells.Item($row, 5).Value = "DOWN"
    $sheet.Cells.Item($row, 6).Value = 0.56
    $sheet.Cells.Item($row, 7).Value = 0.68
    $sheet.Cells.Item($row, 8).Value = "CLOSED"
    $sheet.Cells.Item($row, 9).Value = 21.4286
    $sheet.Cells.Item($row, 10).Value = 0.12
    $sheet.Cells.Item($row, 11).Value = 100.54
    $sheet.Cells.Item($row, 12).Value = 0
    $sheet.Cells.Item($row, 13).Value = 0
    $sheet.Cells.Item($row, 14).Value = 0.6
    $sheet.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item($row, 16).Value = "early_exit"
    $sheet.Cells.Item($row, 17).Value = 0.15
}

# ---------------------------------------------------------------------
# All Trades sheet - append row 52
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

# ---------------------------------------------------------------------
# MarketMaking sheet - append row 52
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
